# Actualizacion desde MV -datos-
# - Corrige B93 (77109 -> 77110)
# - Agrega filas diarias de julio 2021 (filas 126-146)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrige el valor existente en B93
$ws.Cells.Item(93, 2).Value = 77110

# Nuevas filas: Serie, M1 (Circulante), D1+Dv+Ahv
$newData = @(
    @("01-07-2021", 82298, 14593, 67705),
    @("02-07-2021", 79998, 14629, 65369),
    @("05-07-2021", 79661, 14700, 64962),
    @("06-07-2021", 79456, 14696, 64760),
    @("07-07-2021", 79710, 14705, 65005),
    @("08-07-2021", 80116, 14734, 65381),
    @("09-07-2021", 80161, 14776, 65385),
    @("12-07-2021", 80177, 14837, 65341),
    @("13-07-2021", 80022, 14801, 65221),
    @("14-07-2021", 80120, 14822, 65298),
    @("15-07-2021", 79852, 14813, 65040),
    @("19-07-2021", 80872, 14880, 65992),
    @("20-07-2021", 79889, 14868, 65020),
    @("21-07-2021", 78852, 14854, 63998),
    @("22-07-2021", 79563, 14847, 64716),
    @("23-07-2021", 79922, 14846, 65076),
    @("26-07-2021", 80881, 14851, 66030),
    @("27-07-2021", 80394, 14815, 65579),
    @("28-07-2021", 81050, 14777, 66273),
    @("29-07-2021", 82076, 14773, 67303),
    @("30-07-2021", 81330, 14804, 66527)
)

$row = 126
foreach ($item in $newData) {
    $serie = $item[0]
    $m1 = $item[1]
    $circulante = $item[2]
    $cuentas = $item[3]

    # Escribe la fecha como texto (no como fecha serializada) vía formula + pegado de valores
    $ws.Cells.Item($row, 1).Formula = "=""" + $serie + """"
    $ws.Cells.Item($row, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = $m1
    $ws.Cells.Item($row, 3).Value = $circulante
    $ws.Cells.Item($row, 4).Value = $cuentas

    $row = $row + 1
}

$excel.CutCopyMode = 0
